$wb = $excel.ActiveWorkbook

# --- Update the suggested "Andamento" values on sheet acervo_4-1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "Distribuído por exclusão de Ministro"
$ws1.Range("B3").Value = "Distribuído"
$ws1.Range("B4").Value = "Distribuído por prevenção"
$ws1.Range("B5").Value = "Distribuído por prevenção de Turma"

# --- Move the active tab / selection back to acervo_4-1 ---
$ws1.Activate()
$ws1.Range("B6").Select()
